$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1442.1428
$ws.Range("I70").Value = 1499.1666
$ws.Range("J70").Value = 1100
$ws.Range("K70").Value = 4497.4998
$ws.Range("L70").Value = 3300
$ws.Range("M70").Value = -4227.4998
$ws.Range("N70").Value = -3840

$ws.Range("H73").Value = 1442.1428
$ws.Range("I73").Value = 1499.1666
$ws.Range("J73").Value = 1100
$ws.Range("K73").Value = 4497.4998
$ws.Range("L73").Value = 3300
$ws.Range("M73").Value = -3561.4998
$ws.Range("N73").Value = -5172

$ws.Range("H74").Value = 13750.5
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14064

$ws.Range("H77").Value = 13750.5
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70320

$ws.Range("H80").Value = 259.7143
$ws.Range("J80").Value = 236.09091
$ws.Range("L80").Value = 708.27273
$ws.Range("N80").Value = -2704.27273

$ws.Range("H83").Value = 259.7143
$ws.Range("J83").Value = 236.09091
$ws.Range("L83").Value = 2124.81819
$ws.Range("N83").Value = -12108.81819

$ws.Range("H93").Value = 108431.664
$ws.Range("J93").Value = 108431.664
$ws.Range("L93").Value = 108431.664
$ws.Range("N93").Value = -113423.664

$ws.Range("H112").Value = 2129.8948
$ws.Range("J112").Value = 2129.8948
$ws.Range("L112").Value = 6389.6844
$ws.Range("N112").Value = -8605.6844

$ws.Range("H118").Value = 519.44446
$ws.Range("I118").Value = 546.875
$ws.Range("J118").Value = 300
$ws.Range("K118").Value = 1640.625
$ws.Range("L118").Value = 900
$ws.Range("M118").Value = 16.375
$ws.Range("N118").Value = -4214

$ws.Range("H132").Value = 2087.8948
$ws.Range("I132").Value = 1969.7222
$ws.Range("K132").Value = 5909.1666
$ws.Range("M132").Value = -3379.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 31999
$ws.Range("J62").Value = 31999
$ws.Range("L62").Value = 31999
$ws.Range("N62").Value = -33247

$ws.Range("H65").Value = 31999
$ws.Range("J65").Value = 31999
$ws.Range("L65").Value = 95997
$ws.Range("N65").Value = -102237

$ws.Range("H94").Value = 69969
$ws.Range("J94").Value = 69969
$ws.Range("L94").Value = 69969
$ws.Range("N94").Value = -71771

$ws.Range("H132").Value = 4892.952
$ws.Range("I132").Value = 4173.1875
$ws.Range("J132").Value = 7196.2
$ws.Range("K132").Value = 12519.5625
$ws.Range("L132").Value = 21588.6
$ws.Range("M132").Value = -9989.5625
$ws.Range("N132").Value = -26648.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 2875
$ws.Range("J88").Value = 2875
$ws.Range("L88").Value = 2875
$ws.Range("N88").Value = -3687

$ws.Range("H91").Value = 2875
$ws.Range("J91").Value = 2875
$ws.Range("L91").Value = 2875
$ws.Range("N91").Value = -5683

$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992

$ws.Range("H94").Value = 975.9032
$ws.Range("I94").Value = 519.0769
$ws.Range("K94").Value = 519.0769
$ws.Range("M94").Value = -68.07690000000002

$ws.Range("H96").Value = 17717.6
$ws.Range("I96").Value = 17717.6
$ws.Range("K96").Value = 17717.6
$ws.Range("M96").Value = -14971.6

$ws.Range("H134").Value = 4386.25
$ws.Range("I134").Value = 2918.5
$ws.Range("K134").Value = 8755.5
$ws.Range("M134").Value = -6220.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 29999
$ws.Range("J23").Value = 29999
$ws.Range("L23").Value = 29999
$ws.Range("N23").Value = -30479

$ws.Range("H27").Value = 29999
$ws.Range("J27").Value = 29999
$ws.Range("L27").Value = 29999
$ws.Range("N27").Value = -30383

$ws.Range("H31").Value = 3288.25
$ws.Range("I31").Value = 2115.647
$ws.Range("K31").Value = 2115.647
$ws.Range("M31").Value = -1820.647

$ws.Range("H34").Value = 3288.25
$ws.Range("I34").Value = 2115.647
$ws.Range("K34").Value = 2115.647
$ws.Range("M34").Value = -1913.647

$ws.Range("H95").Value = 40898.875
$ws.Range("J95").Value = 40898.875
$ws.Range("L95").Value = 40898.875
$ws.Range("N95").Value = -46390.875

$ws.Range("H107").Value = 676.5263
$ws.Range("I107").Value = 589
$ws.Range("K107").Value = 589
$ws.Range("M107").Value = 1331

$ws.Range("H123").Value = 118288
$ws.Range("J123").Value = 118288
$ws.Range("L123").Value = 118288
$ws.Range("N123").Value = -128088

$ws.Range("H125").Value = 74992
$ws.Range("J125").Value = 74992
$ws.Range("L125").Value = 74992
$ws.Range("N125").Value = -79912

$ws.Range("H134").Value = 3128.9375
$ws.Range("I134").Value = 3181.889
$ws.Range("K134").Value = 9545.667000000001
$ws.Range("M134").Value = -7010.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14077.72
$ws.Range("I70").Value = 17573.025
$ws.Range("K70").Value = 17573.025
$ws.Range("M70").Value = -17303.025

$ws.Range("H73").Value = 14077.72
$ws.Range("I73").Value = 17573.025
$ws.Range("K73").Value = 17573.025
$ws.Range("M73").Value = -16637.025

$ws.Range("H95").Value = 29914.666
$ws.Range("J95").Value = 29914.666
$ws.Range("L95").Value = 29914.666
$ws.Range("N95").Value = -35406.666

$ws.Range("H101").Value = 28842.166
$ws.Range("J101").Value = 28842.166
$ws.Range("L101").Value = 28842.166
$ws.Range("N101").Value = -35332.166

$ws.Range("H102").Value = 3847.1667
$ws.Range("I102").Value = 3589.2856
$ws.Range("K102").Value = 3589.2856
$ws.Range("M102").Value = -1967.2856

$ws.Range("H126").Value = 4290.0835
$ws.Range("I126").Value = 3326
$ws.Range("J126").Value = 5639.8
$ws.Range("K126").Value = 9978
$ws.Range("L126").Value = 16919.4
$ws.Range("M126").Value = -7508
$ws.Range("N126").Value = -21859.4

$ws.Range("H132").Value = 3719.6875
$ws.Range("I132").Value = 3687.1428
$ws.Range("K132").Value = 11061.4284
$ws.Range("M132").Value = -8531.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3155
$ws.Range("I22").Value = 496.16666
$ws.Range("J22").Value = 4927.5557
$ws.Range("K22").Value = 496.16666
$ws.Range("L22").Value = 4927.5557
$ws.Range("M22").Value = -201.16666
$ws.Range("N22").Value = -5517.5557

$ws.Range("H27").Value = 3155
$ws.Range("I27").Value = 496.16666
$ws.Range("J27").Value = 4927.5557
$ws.Range("K27").Value = 496.16666
$ws.Range("L27").Value = 4927.5557
$ws.Range("M27").Value = -389.16666
$ws.Range("N27").Value = -5141.5557

$ws.Range("H61").Value = 3909.1333
$ws.Range("I61").Value = 4453.75
$ws.Range("J61").Value = 1730.6666
$ws.Range("K61").Value = 4453.75
$ws.Range("L61").Value = 1730.6666
$ws.Range("M61").Value = -4251.75
$ws.Range("N61").Value = -2134.6666

$ws.Range("H82").Value = 2320.1538
$ws.Range("I82").Value = 2017.6364
$ws.Range("J82").Value = 3984
$ws.Range("K82").Value = 2017.6364
$ws.Range("L82").Value = 3984
$ws.Range("M82").Value = -1656.6364
$ws.Range("N82").Value = -4706

$ws.Range("H85").Value = 2320.1538
$ws.Range("I85").Value = 2017.6364
$ws.Range("J85").Value = 3984
$ws.Range("K85").Value = 2017.6364
$ws.Range("L85").Value = 3984
$ws.Range("M85").Value = -769.6364000000001
$ws.Range("N85").Value = -6480

$ws.Range("H94").Value = 47582
$ws.Range("J94").Value = 47582
$ws.Range("L94").Value = 47582
$ws.Range("N94").Value = -48934

$ws.Range("H113").Value = 3909.1333
$ws.Range("I113").Value = 4453.75
$ws.Range("J113").Value = 1730.6666
$ws.Range("K113").Value = 4453.75
$ws.Range("L113").Value = 1730.6666
$ws.Range("M113").Value = -2283.75
$ws.Range("N113").Value = -6070.6666

$ws.Range("H132").Value = 3297.889
$ws.Range("I132").Value = 2698
$ws.Range("J132").Value = 4857.6
$ws.Range("K132").Value = 8094
$ws.Range("L132").Value = 14572.8
$ws.Range("M132").Value = -5564
$ws.Range("N132").Value = -19632.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 762.8333
$ws.Range("I100").Value = 415.4
$ws.Range("K100").Value = 830.8
$ws.Range("M100").Value = -289.8

$ws.Range("H104").Value = 215000
$ws.Range("J104").Value = 215000
$ws.Range("L104").Value = 215000
$ws.Range("N104").Value = -221988

$ws.Range("H122").Value = 13238371
$ws.Range("I122").Value = 9618470
$ws.Range("J122").Value = 25003048
$ws.Range("K122").Value = 28855410
$ws.Range("L122").Value = 75009144
$ws.Range("M122").Value = -28852960
$ws.Range("N122").Value = -75014044

$ws.Range("H132").Value = 1305.6154
$ws.Range("I132").Value = 1251.909
$ws.Range("J132").Value = 1601
$ws.Range("K132").Value = 3755.727
$ws.Range("L132").Value = 4803
$ws.Range("M132").Value = -1225.727
$ws.Range("N132").Value = -9863
